# Update "想去人数" (want-to-go count) figures for the two sheets that
# carry the full event listing ("展览" and "全部类型"). Each value is
# incremented by 1, matching the refreshed scrape output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 339
    $ws.Range("F3").Value = 241
    $ws.Range("F5").Value = 288
}
